$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 11691
$ws.Range("B2").Value = 1280
$ws.Range("C2").Value = 300
$ws.Range("D2").Value = "Romain"
$ws.Range("E2").Value = "COUPPE"
$ws.Range("F2").Value = "Romain.Couppe@fr.toyota-industries.eu"
$ws.Range("G2").Value = "Toyota Material Handling France S.A.S."
$ws.Range("H2").Value = "LITHIUM-ION TMHMS & TMHMI"
$ws.Range("I2").Value = "LITHIUM-ION TMHMS & TMHMI"
$ws.Range("J2").Value = "1404-T2-TE-61"
$ws.Range("K2").Value = 719
$ws.Range("L2").Value = "CARQUEFOU 2025 - LITHIUM-ION TMHMS & TMHMI"
$ws.Range("M2").Value = 6
$ws.Range("N2").Value = 2012
$ws.Range("O2").Value = "2025-09-03 13:30"
$ws.Range("P2").Value = "2025-09-04 17:30"
$ws.Range("Q2").Value = 12

# Row 3
$ws.Range("A3").Value = 11628
$ws.Range("B3").Value = 1280
$ws.Range("C3").Value = 300
$ws.Range("D3").Value = "Michel"
$ws.Range("E3").Value = "LAUTRIDOU"
$ws.Range("F3").Value = "Michel.Lautridou@fr.toyota-industries.eu"
$ws.Range("G3").Value = "Toyota Material Handling France S.A.S."
$ws.Range("H3").Value = "LITHIUM-ION TMHMS & TMHMI"
$ws.Range("I3").Value = "LITHIUM-ION TMHMS & TMHMI"
$ws.Range("J3").Value = "1404-T2-TE-61"
$ws.Range("K3").Value = 719
$ws.Range("L3").Value = "CARQUEFOU 2025 - LITHIUM-ION TMHMS & TMHMI"
$ws.Range("M3").Value = 6
$ws.Range("N3").Value = 2012
$ws.Range("O3").Value = "2025-09-03 13:30"
$ws.Range("P3").Value = "2025-09-04 17:30"
$ws.Range("Q3").Value = 12

# Row 4
$ws.Range("A4").Value = 11786
$ws.Range("B4").Value = 1280
$ws.Range("C4").Value = 300
$ws.Range("D4").Value = "Stéphane"
$ws.Range("E4").Value = "VILLETTE"
$ws.Range("F4").Value = "Stephane.VILLETTE@fr.toyota-industries.eu"
$ws.Range("G4").Value = "Toyota Material Handling France S.A.S."
$ws.Range("H4").Value = "LITHIUM-ION TMHMS & TMHMI"
$ws.Range("I4").Value = "LITHIUM-ION TMHMS & TMHMI"
$ws.Range("J4").Value = "1404-T2-TE-61"
$ws.Range("K4").Value = 719
$ws.Range("L4").Value = "CARQUEFOU 2025 - LITHIUM-ION TMHMS & TMHMI"
$ws.Range("M4").Value = 6
$ws.Range("N4").Value = 2012
$ws.Range("O4").Value = "2025-09-03 13:30"
$ws.Range("P4").Value = "2025-09-04 17:30"
$ws.Range("Q4").Value = 12

# Row 5
$ws.Range("A5").Value = 11929
$ws.Range("B5").Value = 1280
$ws.Range("C5").Value = 300
$ws.Range("D5").Value = "Hervé"
$ws.Range("E5").Value = "GUION"
$ws.Range("F5").Value = "Herve.GUION@fr.toyota-industries.eu"
$ws.Range("G5").Value = "Toyota Material Handling France S.A.S."
$ws.Range("H5").Value = "LITHIUM-ION TMHMS & TMHMI"
$ws.Range("I5").Value = "LITHIUM-ION TMHMS & TMHMI"
$ws.Range("J5").Value = "1404-T2-TE-61"
$ws.Range("K5").Value = 719
$ws.Range("L5").Value = "CARQUEFOU 2025 - LITHIUM-ION TMHMS & TMHMI"
$ws.Range("M5").Value = 6
$ws.Range("N5").Value = 2012
$ws.Range("O5").Value = "2025-09-03 13:30"
$ws.Range("P5").Value = "2025-09-04 17:30"
$ws.Range("Q5").Value = 12

# Row 6
$ws.Range("A6").Value = 11712
$ws.Range("B6").Value = 1280
$ws.Range("C6").Value = 300
$ws.Range("D6").Value = "Guillaume"
$ws.Range("E6").Value = "TREBUTIEN"
$ws.Range("F6").Value = "Guillaume.TREBUTIEN@fr.toyota-industries.eu"
$ws.Range("G6").Value = "Toyota Material Handling France S.A.S."
$ws.Range("H6").Value = "LITHIUM-ION TMHMS & TMHMI"
$ws.Range("I6").Value = "LITHIUM-ION TMHMS & TMHMI"
$ws.Range("J6").Value = "1404-T2-TE-61"
$ws.Range("K6").Value = 719
$ws.Range("L6").Value = "CARQUEFOU 2025 - LITHIUM-ION TMHMS & TMHMI"
$ws.Range("M6").Value = 6
$ws.Range("N6").Value = 2012
$ws.Range("O6").Value = "2025-09-03 13:30"
$ws.Range("P6").Value = "2025-09-04 17:30"
$ws.Range("Q6").Value = 12

# Row 7
$ws.Range("A7").Value = 15809
$ws.Range("B7").Value = 1280
$ws.Range("C7").Value = 300
$ws.Range("D7").Value = "Jean-Jacques"
$ws.Range("E7").Value = "MUGABE"
$ws.Range("F7").Value = "Jean-Jacques.MUGABE@fr.toyota-industries.eu"
$ws.Range("G7").Value = "Toyota Material Handling France S.A.S."
$ws.Range("H7").Value = "LITHIUM-ION TMHMS & TMHMI"
$ws.Range("I7").Value = "LITHIUM-ION TMHMS & TMHMI"
$ws.Range("J7").Value = "1404-T2-TE-61"
$ws.Range("K7").Value = 719
$ws.Range("L7").Value = "CARQUEFOU 2025 - LITHIUM-ION TMHMS & TMHMI"
$ws.Range("M7").Value = 6
$ws.Range("N7").Value = 2012
$ws.Range("O7").Value = "2025-09-03 13:30"
$ws.Range("P7").Value = "2025-09-04 17:30"
$ws.Range("Q7").Value = 12
